# DISARM_TAGGING_WORKBOOK.xlsx edit:
#  - Convert the "Table_1" ListObject on "DISARM Red with IDs" into a plain
#    range (removing the table definition) while preserving the banded
#    row colors the table style was providing.
#  - Increase the header row height on both sheets.
#  - Update the saved selection on both sheets to A1:C1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Convert the table to a normal range, dropping the table/ListObject
# definition (xl/tables/table1.xml) but keeping the cell data & most
# formatting intact.
$lo = $ws1.ListObjects.Item(1)
$lo.Unlist()

# The table style used banded rows (white / light gray). Once the table
# is removed those stripes are no longer rendered automatically, so bake
# them into the cell fills directly for the data rows (3-33), matching
# what Excel itself does on "Convert to Range".
for ($r = 3; $r -le 33; $r++) {
    $rowRange = $ws1.Range("A" + $r + ":P" + $r)
    if (($r % 2) -eq 1) {
        # odd data rows -> first row stripe (white)
        $rowRange.Interior.Color = 16777215
    } else {
        # even data rows -> second row stripe (light gray)
        $rowRange.Interior.Color = 15987699
    }
}

# Taller header rows on both sheets.
$ws1.Rows.Item(1).RowHeight = 46.5
$ws2.Rows.Item(1).RowHeight = 37.5

# Update the stored selection on each sheet to A1:C1. Select sheet2 first
# so that sheet1 ends up as the active/selected sheet, matching the
# original workbook (tabSelected on sheet1, gridlines hidden on sheet2).
$ws2.Range("A1:C1").Select()
$ws1.Range("A1:C1").Select()
